# Update the "Summary" sheet: rename contract-related metrics to the
# surviving "Strategy ..." labels, drop the now-redundant contract/strategy
# duplicate rows, and drop the embedded open-positions mini table (row 12).
$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A2").Value = "Strategy Trades"
$summary.Range("A3").Value = "Strategy Win Rate"
$summary.Range("A4").Value = "Strategy Total PnL"

# Remove rows 5-9 (Contract avg hold, Strategy trades/win rate/total PnL/avg
# hold); row 10 ("Verdict") shifts up to become row 5.
$summary.Range("A5:A9").EntireRow.Delete()

# Remove the embedded table header (now shifted from row 12 to row 7).
$summary.Range("A7").EntireRow.Delete()

# Update the "Open Positions" sheet: remove its only row (the header),
# leaving an empty sheet.
$openPositions = $wb.Worksheets.Item("Open Positions")
$openPositions.Range("A1").EntireRow.Delete()
